{"js": "// Update the four r\u00e9sum\u00e9 entries (work-experience job titles + the\n// in-progress education line) so each gains its organization/context\n// prefix and a spaced en-dash in the date range, per the commit diff.\nconst replacements = [\n  {\n    find: \"\uc560\ub2c8\uba54\uc774\uc158 \ub514\uc790\uc774\ub108(2021\ub144 1\uc6d4~\ud604\uc7ac)\",\n    replace: \"Spark \uc560\ub2c8\uba54\uc774\uc158: \uc560\ub2c8\uba54\uc774\uc158 \ub514\uc790\uc774\ub108(2021\ub144 1\uc6d4 - \ud604\uc7ac)\",\n  },\n  {\n    find: \"\uc560\ub2c8\uba54\uc774\uc158 \ub514\uc790\uc774\ub108(2018\ub144 6\uc6d4~2020\ub144 12\uc6d4)\",\n    replace: \"Pixel Studio: \uc560\ub2c8\uba54\uc774\uc158 \ub514\uc790\uc774\ub108(2018\ub144 6\uc6d4 - 2020\ub144 12\uc6d4)\",\n  },\n  {\n    find: \"\ubcf4\uc870 \uc560\ub2c8\uba54\uc774\uc158 \ub514\uc790\uc774\ub108(2016\ub144 9\uc6d4~2018\ub144 5\uc6d4)\",\n    replace: \"\ud50c\ub798\uc2dc \uc560\ub2c8\uba54\uc774\uc158: \uc8fc\ub2c8\uc5b4 \uc560\ub2c8\uba54\uc774\uc158 \ub514\uc790\uc774\ub108(2016\ub144 9\uc6d4 - 2018\ub144 5\uc6d4)\",\n  },\n  {\n    find: \"\uc608\uc220\ud559\ubd80 \uc560\ub2c8\uba54\uc774\uc158\uacfc \uc11d\uc0ac \ud559\uc704 \ucde8\ub4dd \uc608\uc815(\uc878\uc5c5 \uc608\uc815\uc77c:\",\n    replace: \"\uc560\ub2c8\uba54\uc774\uc158\uc758 \uc608\uc220 \ub9c8\uc2a4\ud130, \uc608\uc0c1 \uc878\uc5c5: 2025\ub144 12\uc6d4\",\n  },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${find}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the four r\u00e9sum\u00e9 entries (work-experience job titles + the\n# in-progress education line) so each gains its organization/context\n# prefix and a spaced en-dash in the date range, per the commit diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"\uc560\ub2c8\uba54\uc774\uc158 \ub514\uc790\uc774\ub108(2021\ub144 1\uc6d4~\ud604\uc7ac)\"; Replace = \"Spark \uc560\ub2c8\uba54\uc774\uc158: \uc560\ub2c8\uba54\uc774\uc158 \ub514\uc790\uc774\ub108(2021\ub144 1\uc6d4 - \ud604\uc7ac)\" },\n    @{ Find = \"\uc560\ub2c8\uba54\uc774\uc158 \ub514\uc790\uc774\ub108(2018\ub144 6\uc6d4~2020\ub144 12\uc6d4)\"; Replace = \"Pixel Studio: \uc560\ub2c8\uba54\uc774\uc158 \ub514\uc790\uc774\ub108(2018\ub144 6\uc6d4 - 2020\ub144 12\uc6d4)\" },\n    @{ Find = \"\ubcf4\uc870 \uc560\ub2c8\uba54\uc774\uc158 \ub514\uc790\uc774\ub108(2016\ub144 9\uc6d4~2018\ub144 5\uc6d4)\"; Replace = \"\ud50c\ub798\uc2dc \uc560\ub2c8\uba54\uc774\uc158: \uc8fc\ub2c8\uc5b4 \uc560\ub2c8\uba54\uc774\uc158 \ub514\uc790\uc774\ub108(2016\ub144 9\uc6d4 - 2018\ub144 5\uc6d4)\" },\n    @{ Find = \"\uc608\uc220\ud559\ubd80 \uc560\ub2c8\uba54\uc774\uc158\uacfc \uc11d\uc0ac \ud559\uc704 \ucde8\ub4dd \uc608\uc815(\uc878\uc5c5 \uc608\uc815\uc77c:\"; Replace = \"\uc560\ub2c8\uba54\uc774\uc158\uc758 \uc608\uc220 \ub9c8\uc2a4\ud130, \uc608\uc0c1 \uc878\uc5c5: 2025\ub144 12\uc6d4\" }\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $find.Execute($r.Find, $false, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2) | Out-Null\n}\n"}
